{"js": "// Update the date line and the twenty-five \"three-digit \u00f7 one-digit\"\n// answer cells to the new day's generated values.\nconst replacements = [\n  [\"2025-02-22 Saturday\", \"2025-02-23 Sunday\"],\n  [\"220\u00f74=55, 0\", \"945\u00f76=157, 3\"],\n  [\"771\u00f73=257, 0\", \"531\u00f75=106, 1\"],\n  [\"332\u00f74=83, 0\", \"725\u00f75=145, 0\"],\n  [\"331\u00f75=66, 1\", \"104\u00f72=52, 0\"],\n  [\"391\u00f76=65, 1\", \"536\u00f75=107, 1\"],\n  [\"236\u00f74=59, 0\", \"160\u00f73=53, 1\"],\n  [\"247\u00f73=82, 1\", \"904\u00f74=226, 0\"],\n  [\"533\u00f76=88, 5\", \"909\u00f74=227, 1\"],\n  [\"527\u00f78=65, 7\", \"575\u00f78=71, 7\"],\n  [\"525\u00f79=58, 3\", \"768\u00f78=96, 0\"],\n  [\"194\u00f76=32, 2\", \"705\u00f79=78, 3\"],\n  [\"992\u00f79=110, 2\", \"624\u00f74=156, 0\"],\n  [\"737\u00f76=122, 5\", \"570\u00f79=63, 3\"],\n  [\"978\u00f73=326, 0\", \"620\u00f76=103, 2\"],\n  [\"907\u00f76=151, 1\", \"878\u00f79=97, 5\"],\n  [\"758\u00f74=189, 2\", \"433\u00f78=54, 1\"],\n  [\"262\u00f78=32, 6\", \"701\u00f76=116, 5\"],\n  [\"442\u00f77=63, 1\", \"941\u00f73=313, 2\"],\n  [\"144\u00f75=28, 4\", \"311\u00f73=103, 2\"],\n  [\"513\u00f74=128, 1\", \"554\u00f79=61, 5\"],\n  [\"543\u00f75=108, 3\", \"481\u00f77=68, 5\"],\n  [\"229\u00f76=38, 1\", \"948\u00f79=105, 3\"],\n  [\"359\u00f73=119, 2\", \"418\u00f74=104, 2\"],\n  [\"885\u00f76=147, 3\", \"863\u00f79=95, 8\"],\n  [\"578\u00f77=82, 4\", \"282\u00f73=94, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty-five \"three-digit \u00f7 one-digit\"\n# answer cells to the new day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-02-22 Saturday\", \"2025-02-23 Sunday\"),\n  @(\"220\u00f74=55, 0\", \"945\u00f76=157, 3\"),\n  @(\"771\u00f73=257, 0\", \"531\u00f75=106, 1\"),\n  @(\"332\u00f74=83, 0\", \"725\u00f75=145, 0\"),\n  @(\"331\u00f75=66, 1\", \"104\u00f72=52, 0\"),\n  @(\"391\u00f76=65, 1\", \"536\u00f75=107, 1\"),\n  @(\"236\u00f74=59, 0\", \"160\u00f73=53, 1\"),\n  @(\"247\u00f73=82, 1\", \"904\u00f74=226, 0\"),\n  @(\"533\u00f76=88, 5\", \"909\u00f74=227, 1\"),\n  @(\"527\u00f78=65, 7\", \"575\u00f78=71, 7\"),\n  @(\"525\u00f79=58, 3\", \"768\u00f78=96, 0\"),\n  @(\"194\u00f76=32, 2\", \"705\u00f79=78, 3\"),\n  @(\"992\u00f79=110, 2\", \"624\u00f74=156, 0\"),\n  @(\"737\u00f76=122, 5\", \"570\u00f79=63, 3\"),\n  @(\"978\u00f73=326, 0\", \"620\u00f76=103, 2\"),\n  @(\"907\u00f76=151, 1\", \"878\u00f79=97, 5\"),\n  @(\"758\u00f74=189, 2\", \"433\u00f78=54, 1\"),\n  @(\"262\u00f78=32, 6\", \"701\u00f76=116, 5\"),\n  @(\"442\u00f77=63, 1\", \"941\u00f73=313, 2\"),\n  @(\"144\u00f75=28, 4\", \"311\u00f73=103, 2\"),\n  @(\"513\u00f74=128, 1\", \"554\u00f79=61, 5\"),\n  @(\"543\u00f75=108, 3\", \"481\u00f77=68, 5\"),\n  @(\"229\u00f76=38, 1\", \"948\u00f79=105, 3\"),\n  @(\"359\u00f73=119, 2\", \"418\u00f74=104, 2\"),\n  @(\"885\u00f76=147, 3\", \"863\u00f79=95, 8\"),\n  @(\"578\u00f77=82, 4\", \"282\u00f73=94, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $null = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
